$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing data rows (blocks 5-8 are gone; the
# block list was re-chunked from 8x20 rows into 4x24 rows).
$ws.Range("A6:B9").EntireRow.Delete()

# Strip the header's bold/shaded formatting and the alternating
# shaded-row formatting from the remaining rows; only B5 keeps a
# (text) number format further down.
$ws.Rows.Item(1).ClearFormats()
$ws.Range("A2:B3").ClearFormats()

# Rewrite the table contents for the new 24-row-per-block chunking.
$ws.Range("A1").Value = "block_number_loop"
$ws.Range("B1").Value = "list_of_rows"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95"

# B5 carries the same text/"@" number format that used to live on the
# bottom two rows of the old table.
$ws.Range("B5").NumberFormat = "@"

# Restore the cursor to where the author left it in the saved file.
$null = $ws.Range("B8").Select()
